# Append two new data rows to the existing table (rows 4 and 5), reusing
# the same shared-string values already used by rows 2/3 ("tcm:2-64-32",
# "tcm:2-18"), then move the active selection to C7 to match the saved
# workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: # = 2, TemplateId = tcm:2-64-32, ItemId = tcm:2-18, Ignore = TRUE
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "tcm:2-64-32"
$ws.Range("C4").Value = "tcm:2-18"
$ws.Range("D4").Value = $true

# Row 5: # = 3, TemplateId = tcm:2-64-32, ItemId = tcm:2-18, Ignore = FALSE
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "tcm:2-64-32"
$ws.Range("C5").Value = "tcm:2-18"
$ws.Range("D5").Value = $false

# Match the saved selection state (C7) from the authored workbook.
[void]$ws.Range("C7").Select()
